# "Some fixes + wallet controller added"
#
# Add a "CreatorUserID" column to the Suggestions sheet and make that
# sheet the active tab (it was previously the Issues sheet).

$wb = $excel.ActiveWorkbook

$wsSuggestions = $wb.Worksheets.Item("Suggestions")

# New column E: CreatorUserID (header reuses the existing "CreatorUserID"
# shared string already used elsewhere in the workbook, e.g. on the
# WalletTransactions sheet's WalletID header style).
$wsSuggestions.Range("E1").Value = "CreatorUserID"
$wsSuggestions.Range("E1").Font.Bold = $true
$wsSuggestions.Range("E1").HorizontalAlignment = -4108   # xlCenter

$wsSuggestions.Range("E2").Value = 4
$wsSuggestions.Range("E3").Value = 5

# Best-fit the new column to its header text ("CreatorUserID").
$wsSuggestions.Columns.Item(5).ColumnWidth = 12.6

# Suggestions becomes the active/selected sheet (previously Issues was).
$wsSuggestions.Activate()
